# 1. Merge the three runs describing the Q1 heading into a single run of text.
#    "n the " + "relations between categorical and numerical" + " variables"
#    -> "n the relations between categorical and numerical variables"
$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "n the relations between categorical and numerical variables",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "n the relations between categorical and numerical variables",
    2) | Out-Null

# 2. Fix the typo ANOVE -> ANOVA within "No test performed, ANOVE might be used"
$d.Content.Find.Execute(
    "No test performed, ANOVE might be used",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "No test performed, ANOVA might be used",
    2) | Out-Null

# 3. Add three new paragraphs after the last "WGHT" bullet item (end of document,
#    just before the section break): an empty paragraph, a paragraph with the
#    "Notes to add on the redaction of Q1: " text, and a trailing empty paragraph.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$insertRange = $lastPara.Range
$insertRange.Collapse(0)  # wdCollapseEnd

$insertRange.InsertParagraphAfter()
$insertRange.Collapse(0)
$insertRange.InsertParagraphAfter()
$insertRange.Collapse(0)
$insertRange.InsertParagraphAfter()
$insertRange.Collapse(0)

# Now set the text of the middle of the three newly inserted paragraphs.
$paraCount = $d.Paragraphs.Count
$notesPara = $d.Paragraphs($paraCount - 1)
$notesRange = $notesPara.Range
$notesRange.Collapse(0)
$notesRange.InsertBefore("Notes to add on the redaction of Q1: ")
